# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G, header "K" in row 1) previously held a
# "Strike#" count. It is regenerated here against the real strikeout (K)
# totals for each outing, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 1
    3  = 1
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 3
    10 = 4
    11 = 1
    12 = 0
    13 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 2
    19 = 3
    20 = 0
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 3
    29 = 2
    30 = 3
    31 = 1
    32 = 2
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 1
    39 = 2
    40 = 1
    41 = 3
    42 = 0
    43 = 1
    44 = 2
    45 = 2
    46 = 2
    47 = 0
    48 = 1
    49 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
